# Regenerate save_data: recompute the K column (column G) values using
# strikeout counts ("K") instead of the previous Strike# metric, then
# rewrite the recalculated stats back into the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values keyed by worksheet row number (row 1 is the header row).
$kValues = @{
    2 = 5
    3 = 6
    4 = 6
    5 = 5
    6 = 6
    7 = 10
    8 = 9
    9 = 4
    10 = 6
    11 = 10
    12 = 6
    13 = 5
    14 = 1
    15 = 5
    16 = 1
    17 = 3
    18 = 3
    19 = 2
    20 = 2
    21 = 4
    22 = 2
    23 = 1
    24 = 8
    25 = 4
    26 = 2
    27 = 3
    28 = 0
    29 = 2
    30 = 1
    31 = 0
    32 = 2
    33 = 1
    34 = 0
    35 = 1
    36 = 0
    37 = 0
    38 = 2
    39 = 0
    40 = 1
    41 = 2
    42 = 1
    43 = 1
    44 = 3
    45 = 2
    46 = 1
    47 = 2
    48 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
